$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 06:42"

# Row 15 - Pakistan
$ws.Range("B15").Value = 263496
$ws.Range("C15").Value = 1580
$ws.Range("D15").Value = 204276
$ws.Range("E15").Value = 53652
$ws.Range("G15").Value = 46
$ws.Range("H15").Value = 5568

# Row 36 - Belgica
$ws.Range("B36").Value = 63706
$ws.Range("C36").Value = 207
$ws.Range("E36").Value = 36617

# Row 89 - Haiti
$ws.Range("B89").Value = 7053
$ws.Range("C89").Value = 78
$ws.Range("D89").Value = 3877
$ws.Range("E89").Value = 3030

# Row 156 - Malta
$ws.Range("B156").Value = 676
$ws.Range("E156").Value = 5

# Row 170 - Mongolia
$ws.Range("D170").Value = 213
$ws.Range("E170").Value = 74
